# Scheduled-runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the per-job "Profits" sheets. Generated from the upstream diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 9767
$ws.Range("J88").Value = 9767
$ws.Range("L88").Value = 9767
$ws.Range("N88").Value = -10579

$ws.Range("H91").Value = 9767
$ws.Range("J91").Value = 9767
$ws.Range("L91").Value = 9767
$ws.Range("N91").Value = -12575

$ws.Range("H99").Value = 1564.6666
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -3002

$ws.Range("H113").Value = 5220.3335
$ws.Range("I113").Value = 2830.5
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2830.5
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 423.5
$ws.Range("N113").Value = -16508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 1475
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1475
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1475
$ws.Range("N29").Value = -2091
$ws.Range("M29").ClearContents()

$ws.Range("H63").Value = 2951.25
$ws.Range("I63").Value = 1902.5
$ws.Range("K63").Value = 1902.5
$ws.Range("M63").Value = -1216.5

$ws.Range("H66").Value = 2951.25
$ws.Range("I66").Value = 1902.5
$ws.Range("K66").Value = 9512.5
$ws.Range("M66").Value = -6080.5

$ws.Range("H97").Value = 1144.7778
$ws.Range("I97").Value = 1009.6667
$ws.Range("J97").Value = 1415
$ws.Range("K97").Value = 1009.6667
$ws.Range("L97").Value = 1415
$ws.Range("M97").Value = -513.6667
$ws.Range("N97").Value = -2407

$ws.Range("H102").Value = 1573.125
$ws.Range("I102").Value = 1573.125
$ws.Range("K102").Value = 1573.125
$ws.Range("M102").Value = 48.875

$ws.Range("H110").Value = 3530
$ws.Range("I110").Value = 2328.2
$ws.Range("J110").Value = 3992.2307
$ws.Range("K110").Value = 2328.2
$ws.Range("L110").Value = 3992.2307
$ws.Range("M110").Value = -283.1999999999998
$ws.Range("N110").Value = -8082.2307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1803.6296
$ws.Range("I94").Value = 1812.2084
$ws.Range("K94").Value = 1812.2084
$ws.Range("M94").Value = -1361.2084

$ws.Range("H99").Value = 3833.3333
$ws.Range("I99").Value = 3833.3333
$ws.Range("K99").Value = 3833.3333
$ws.Range("M99").Value = -2335.3333

$ws.Range("H105").Value = 2980.4
$ws.Range("I105").Value = 2804.8333
$ws.Range("K105").Value = 2804.8333
$ws.Range("M105").Value = -1057.8333

$ws.Range("H134").Value = 6859.2666
$ws.Range("I134").Value = 6859.2666
$ws.Range("K134").Value = 20577.7998
$ws.Range("M134").Value = -18042.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8199.857
$ws.Range("I16").Value = 3466.6667
$ws.Range("J16").Value = 11749.75
$ws.Range("K16").Value = 3466.6667
$ws.Range("L16").Value = 11749.75
$ws.Range("M16").Value = -3179.6667
$ws.Range("N16").Value = -12323.75

$ws.Range("H31").Value = 3521.077
$ws.Range("I31").Value = 2594.9
$ws.Range("J31").Value = 6608.3335
$ws.Range("K31").Value = 2594.9
$ws.Range("L31").Value = 6608.3335
$ws.Range("M31").Value = -2299.9
$ws.Range("N31").Value = -7198.3335

$ws.Range("H34").Value = 3521.077
$ws.Range("I34").Value = 2594.9
$ws.Range("J34").Value = 6608.3335
$ws.Range("K34").Value = 2594.9
$ws.Range("L34").Value = 6608.3335
$ws.Range("M34").Value = -2392.9
$ws.Range("N34").Value = -7012.3335

$ws.Range("H113").Value = 8199.857
$ws.Range("I113").Value = 3466.6667
$ws.Range("J113").Value = 11749.75
$ws.Range("K113").Value = 3466.6667
$ws.Range("L113").Value = 11749.75
$ws.Range("M113").Value = -1296.6667
$ws.Range("N113").Value = -16089.75

$ws.Range("H122").Value = 2465.8
$ws.Range("I122").Value = 2571.5881
$ws.Range("K122").Value = 7714.7643
$ws.Range("M122").Value = -5264.7643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 521.1579
$ws.Range("I5").Value = 475.3125
$ws.Range("K5").Value = 1425.9375
$ws.Range("M5").Value = -1313.9375

$ws.Range("H135").Value = 521.1579
$ws.Range("I135").Value = 475.3125
$ws.Range("K135").Value = 4277.8125
$ws.Range("M135").Value = -1742.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1398.9474
$ws.Range("I102").Value = 1393.3334
$ws.Range("K102").Value = 1393.3334
$ws.Range("M102").Value = 228.6666

$ws.Range("H122").Value = 17827.62
$ws.Range("I122").Value = 13716.529
$ws.Range("J122").Value = 35299.75
$ws.Range("K122").Value = 41149.587
$ws.Range("L122").Value = 105899.25
$ws.Range("M122").Value = -38699.587
$ws.Range("N122").Value = -110799.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4185.278
$ws.Range("I40").Value = 3181.2144
$ws.Range("J40").Value = 7699.5
$ws.Range("K40").Value = 3181.2144
$ws.Range("L40").Value = 7699.5
$ws.Range("M40").Value = -3045.2144
$ws.Range("N40").Value = -7971.5

$ws.Range("H46").Value = 1771.5454
$ws.Range("I46").Value = 883.3333
$ws.Range("J46").Value = 2837.4
$ws.Range("K46").Value = 883.3333
$ws.Range("L46").Value = 2837.4
$ws.Range("M46").Value = -695.3333
$ws.Range("N46").Value = -3213.4

$ws.Range("H56").Value = 10051
$ws.Range("I56").Value = 10051
$ws.Range("K56").Value = 10051
$ws.Range("M56").Value = -9360

$ws.Range("H136").Value = 3509.3572
$ws.Range("I136").Value = 3302.3635
$ws.Range("J136").Value = 4268.3335
$ws.Range("K136").Value = 9907.0905
$ws.Range("L136").Value = 12805.0005
$ws.Range("M136").Value = -7357.0905
$ws.Range("N136").Value = -17905.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 35000
$ws.Range("I51").Value = 29000
$ws.Range("J51").Value = 41000
$ws.Range("K51").Value = 29000
$ws.Range("L51").Value = 41000
$ws.Range("M51").Value = -28490
$ws.Range("N51").Value = -42020

$ws.Range("H132").Value = 641.8
$ws.Range("I132").Value = 604.3333
$ws.Range("J132").Value = 698
$ws.Range("K132").Value = 1812.9999
$ws.Range("L132").Value = 2094
$ws.Range("M132").Value = 717.0001
$ws.Range("N132").Value = -7154

$ws.Range("H136").Value = 3163.647
$ws.Range("I136").Value = 3212.1333
$ws.Range("K136").Value = 9636.3999
$ws.Range("M136").Value = -7086.3999
